$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Add a new data row (row 3) with the new sample record
$ws.Range("A3").Value = "rajan"
$ws.Range("B3").Value = 12021986
$ws.Range("C3").Value = "12Kanmvle"
$ws.Range("D3").Value = "palani"
$ws.Range("E3").Value = "chennai"
$ws.Range("F3").Value = 78451269
$ws.Range("G3").Value = 784512
$ws.Range("H3").Value = "mlkdj@gmd.com"
$ws.Range("I3").Value = 124536

# Turn the new email cell into a mailto hyperlink
$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:mlkdj@gmd.com")

# Match the Hyperlink cell style already used by H2, then restore the value
# (Hyperlinks.Add stamps its own formatting, so re-apply the existing look)
$ws.Range("H2").Copy($ws.Range("H3"))
$ws.Range("H3").Value = "mlkdj@gmd.com"

# Update the active selection to reflect where the user left off editing
$ws.Range("H4").Select()
